$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new word entries in the exact order they were originally authored
# (word column first for a batch of rows, then meanings/examples), so the
# resulting shared-string table ordering matches the source edit.
$ws.Cells.Item(27, 1).Value = "correlation"
$ws.Cells.Item(27, 2).Value = "a connection between two ideas, facts etc, especially when one may be the cause of the other"
$ws.Cells.Item(27, 3).Value = "there is a strong correlation between imcome and education"
$ws.Cells.Item(27, 4).Value = "a strong correlation between urban deprivation and poor health"
$ws.Cells.Item(28, 1).Value = "defy"
$ws.Cells.Item(28, 2).Value = "to refuse to obey a law or rule, or refuse to do what someone in authority tells you to do"
$ws.Cells.Item(28, 3).Value = "it's worng to defy the orders of superior officer."
$ws.Cells.Item(28, 4).Value = "Billy defied his mother, and smoked openly in the house."
$ws.Cells.Item(29, 1).Value = "superior"
$ws.Cells.Item(29, 2).Value = " better, more powerful, more effective etc than a similar person or thing, especially one that you are competing against "
$ws.Cells.Item(29, 4).Value = "Fletcher’s superior technique brought him victory."
$ws.Cells.Item(29, 3).Value = "Your computer is far superior to mine."
$ws.Cells.Item(30, 1).Value = "decontaminate"
$ws.Cells.Item(30, 2).Value = "to remove a dangerous substance from somewhere"
$ws.Cells.Item(30, 4).Value = "It may cost over `$5 million to decontaminate the whole site."
$ws.Cells.Item(30, 3).Value = "it also helps decontaminate water and reduce the impact of floods."
$ws.Cells.Item(31, 1).Value = "impact"
$ws.Cells.Item(31, 2).Value = "the effect or influence that an event, situation etc has on someone or something"
$ws.Cells.Item(31, 4).Value = "an international meeting to consider the environmental impacts of global warming"
$ws.Cells.Item(31, 3).Value = "We need to assess the impact on climate change."
$ws.Cells.Item(32, 1).Value = "intention"
$ws.Cells.Item(32, 2).Value = "a plan or desire to do something "
$ws.Cells.Item(32, 4).Value = "I have no intention of retiring just yet."
$ws.Cells.Item(32, 3).Value = "you have the intention permanently to deprive me of the gasoline."
$ws.Cells.Item(33, 1).Value = "lucid"
$ws.Cells.Item(34, 1).Value = "lure"
$ws.Cells.Item(34, 2).Value = "to persuade someone to do something, especially something wrong or dangerous, by making it seem attractive or exciting"
$ws.Cells.Item(34, 3).Value = "cheese is realy good lure for mice."
$ws.Cells.Item(34, 4).Value = "People may be lured into buying tickets by clever advertising."
$ws.Cells.Item(33, 2).Value = "expressed in a way that is clear and easy to understand"
$ws.Cells.Item(33, 3).Value = "after finishing , she became lucid ,recognizedhim , agreed him."
$ws.Cells.Item(33, 4).Value = "You must write in a clear and lucid style."
$ws.Cells.Item(35, 1).Value = "desecrate"
$ws.Cells.Item(36, 1).Value = "invade"
$ws.Cells.Item(37, 1).Value = "devastate"
$ws.Cells.Item(39, 1).Value = "evacuate"
$ws.Cells.Item(40, 1).Value = "stroke"
$ws.Cells.Item(41, 1).Value = "flock"
$ws.Cells.Item(35, 2).Value = "to spoil or damage something holy or respected"
$ws.Cells.Item(35, 4).Value = "Kelly's grave was also desecrated."
$ws.Cells.Item(35, 3).Value = "Most of the Egyptian tombs were desecrated and robbed."
$ws.Cells.Item(36, 2).Value = "to enter a country, town, or area using military force, in order to take control of it"
$ws.Cells.Item(36, 4).Value = "Every summer, the town is invaded by tourists."
$ws.Cells.Item(36, 3).Value = "the invading army desectared this holy place when they camped."
$ws.Cells.Item(37, 3).Value = "each one of bombs can devastate a city."
$ws.Cells.Item(37, 2).Value = "to damage something very badly or completely"
$ws.Cells.Item(37, 4).Value = "The city centre was devastated by the bomb."
$ws.Cells.Item(38, 1).Value = "discreet"
$ws.Cells.Item(38, 2).Value = " careful about what you say or do, so that you do not offend, upset, or embarrass people or tell secrets"
$ws.Cells.Item(38, 3).Value = "she is very discreed in giving her opinion"
$ws.Cells.Item(38, 4).Value = "I stood back at a discreet distance."
$ws.Cells.Item(39, 2).Value = " to send people away from a dangerous place to a safe place"
$ws.Cells.Item(39, 3).Value = "the official ordered the residents to evacuate."
$ws.Cells.Item(39, 4).Value = "During the war he was evacuated to Scotland."
$ws.Cells.Item(40, 2).Value = "if someone has a stroke, an artery (=tube carrying blood) in their brain suddenly bursts or becomes blocked, so that they may die or be unable to use some muscles"
$ws.Cells.Item(40, 3).Value = "they intended to devastate the tower at a stroke."
$ws.Cells.Item(40, 4).Value = "I looked after my father after he had a stroke."
$ws.Cells.Item(41, 2).Value = " a group of sheep, goats, or birds"
$ws.Cells.Item(41, 4).Value = "a flock of small birds"
$ws.Cells.Item(41, 3).Value = "he keeps a flock of sheeps."

# Row heights for the newly added rows (27-41), matching the authored sizing
$ws.Rows.Item(27).RowHeight = 75
$ws.Rows.Item(28).RowHeight = 60
$ws.Rows.Item(29).RowHeight = 75
$ws.Rows.Item(30).RowHeight = 75
$ws.Rows.Item(31).RowHeight = 75
$ws.Rows.Item(32).RowHeight = 75
$ws.Rows.Item(33).RowHeight = 60
$ws.Rows.Item(34).RowHeight = 105
$ws.Rows.Item(35).RowHeight = 60
$ws.Rows.Item(36).RowHeight = 60
$ws.Rows.Item(37).RowHeight = 45
$ws.Rows.Item(38).RowHeight = 75
$ws.Rows.Item(39).RowHeight = 45
$ws.Rows.Item(40).RowHeight = 120
$ws.Rows.Item(41).RowHeight = 30

# Restore the view: scrolled down to show the newly added rows, with E39 selected
$excel.ActiveWindow.ScrollRow = 35
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E39").Select()

